{"js": "// Update the \"Number of Genotypes\" column in the outplant ACER table with\n// the newest NFWF report figures, and refresh the Totals row accordingly.\n//\n// Table layout (row 0 = header):\n//   0: Outplant Date | Site ID | Number of Corals | Latitude | Longitude | Number of Genotypes\n//   1..17: data rows\n//   18: Totals row (Number of Corals total in col 4, Number of Genotypes total in col 5)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New \"Number of Genotypes\" values for data rows 1-17 (in row order).\nconst newGenotypeCounts = [10, 8, 7, 10, 10, 7, 8, 6, 10, 5, 10, 10, 11, 11, 10, 10, 10];\n\nconst genotypeColumnIndex = 5;\n\nfor (let i = 0; i < newGenotypeCounts.length; i++) {\n  const rowIndex = i + 1; // data rows start at 1 (row 0 is the header)\n  const cell = table.getCell(rowIndex, genotypeColumnIndex);\n  cell.value = String(newGenotypeCounts[i]);\n}\n\n// Totals row: recompute the \"Number of Genotypes\" total as the sum of the\n// updated per-row values.\nconst totalsRowIndex = newGenotypeCounts.length + 1; // row 18\nconst total = newGenotypeCounts.reduce((sum, v) => sum + v, 0);\nconst totalsCell = table.getCell(totalsRowIndex, genotypeColumnIndex);\ntotalsCell.value = String(total);\n\nawait context.sync();\n", "ps1": "# Update the \"Number of Genotypes\" column in the outplant ACER table with\n# the newest NFWF report figures, and refresh the Totals row accordingly.\n#\n# Table layout (row 1 = header in Word's 1-based Rows collection):\n#   1: Outplant Date | Site ID | Number of Corals | Latitude | Longitude | Number of Genotypes\n#   2..18: data rows\n#   19: Totals row (Number of Corals total in col 5, Number of Genotypes total in col 6)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New \"Number of Genotypes\" values for data rows, in row order.\n$newVals = @(10, 8, 7, 10, 10, 7, 8, 6, 10, 5, 10, 10, 11, 11, 10, 10, 10)\n$genotypeColumn = 6\n\nfor ($i = 0; $i -lt $newVals.Length; $i++) {\n    $rowIndex = $i + 2\n    $cell = $t.Cell($rowIndex, $genotypeColumn)\n    $rng = $cell.Range\n    $rng.End = $rng.End - 1\n    $rng.Text = [string]$newVals[$i]\n}\n\n# Totals row: recompute the \"Number of Genotypes\" total as the sum of the\n# updated per-row values.\n$total = 0\nforeach ($v in $newVals) { $total += $v }\n\n$totalsRowIndex = $newVals.Length + 2\n$totalsCell = $t.Cell($totalsRowIndex, $genotypeColumn)\n$totalsRange = $totalsCell.Range\n$totalsRange.End = $totalsRange.End - 1\n$totalsRange.Text = [string]$total\n"}
